# Auto-generated script to apply diff changes to Balmung_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1208.6
$ws.Range("I39").Value = 575.55554
$ws.Range("J39").Value = 2158.1667
$ws.Range("K39").Value = 1726.66662
$ws.Range("L39").Value = 6474.500100000001
$ws.Range("M39").Value = -1430.66662
$ws.Range("N39").Value = -7066.500100000001
$ws.Range("H70").Value = 2131
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2131
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6393
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6933
$ws.Range("H73").Value = 2131
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2131
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6393
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -8265
$ws.Range("H111").Value = 739.5
$ws.Range("I111").Value = 654.2857
$ws.Range("K111").Value = 1962.8571
$ws.Range("M111").Value = 1104.1429
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H129").Value = 1801.7273
$ws.Range("I129").Value = 1263.8
$ws.Range("K129").Value = 3791.4
$ws.Range("M129").Value = 1208.6
$ws.Range("H132").Value = 1710.2
$ws.Range("I132").Value = 1720.8182
$ws.Range("K132").Value = 5162.4546
$ws.Range("M132").Value = -2632.4546
$ws.Range("H135").Value = 509.36667
$ws.Range("I135").Value = 350.56522
$ws.Range("K135").Value = 3155.08698
$ws.Range("M135").Value = -620.08698
$ws.Range("H138").Value = 4212.6787
$ws.Range("I138").Value = 5149.75
$ws.Range("J138").Value = 3509.875
$ws.Range("K138").Value = 15449.25
$ws.Range("L138").Value = 10529.625
$ws.Range("M138").Value = -10309.25
$ws.Range("N138").Value = -20809.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 158683.66
$ws.Range("I32").Value = 197269.62
$ws.Range("J32").Value = 7307.923
$ws.Range("K32").Value = 197269.62
$ws.Range("L32").Value = 7307.923
$ws.Range("M32").Value = -196982.62
$ws.Range("N32").Value = -7881.923
$ws.Range("H45").Value = 70044.8
$ws.Range("I45").Value = 86556
$ws.Range("K45").Value = 86556
$ws.Range("M45").Value = -86179
$ws.Range("H61").Value = 54252.027
$ws.Range("I61").Value = 47061.61
$ws.Range("K61").Value = 47061.61
$ws.Range("M61").Value = -46849.61
$ws.Range("H74").Value = 2488.606
$ws.Range("I74").Value = 1415.2273
$ws.Range("J74").Value = 4635.364
$ws.Range("K74").Value = 1415.2273
$ws.Range("L74").Value = 4635.364
$ws.Range("M74").Value = -541.2273
$ws.Range("N74").Value = -6383.364
$ws.Range("H77").Value = 2488.606
$ws.Range("I77").Value = 1415.2273
$ws.Range("J77").Value = 4635.364
$ws.Range("K77").Value = 7076.136500000001
$ws.Range("L77").Value = 23176.82
$ws.Range("M77").Value = -2708.136500000001
$ws.Range("N77").Value = -31912.82
$ws.Range("H102").Value = 2308.1052
$ws.Range("I102").Value = 2308.1052
$ws.Range("K102").Value = 2308.1052
$ws.Range("M102").Value = -686.1052
$ws.Range("H107").Value = 45496
$ws.Range("J107").Value = 45496
$ws.Range("L107").Value = 45496
$ws.Range("N107").Value = -53176
$ws.Range("H132").Value = 3429.6667
$ws.Range("I132").Value = 3077.6924
$ws.Range("K132").Value = 9233.0772
$ws.Range("M132").Value = -6703.0772
$ws.Range("H133").Value = 179990
$ws.Range("J133").Value = 179990
$ws.Range("L133").Value = 179990
$ws.Range("N133").Value = -185050
$ws.Range("H136").Value = 54252.027
$ws.Range("I136").Value = 47061.61
$ws.Range("K136").Value = 141184.83
$ws.Range("M136").Value = -138634.83

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 458.8
$ws.Range("I25").Value = 323.5
$ws.Range("K25").Value = 323.5
$ws.Range("M25").Value = -88.5
$ws.Range("H64").Value = 701.2222
$ws.Range("I64").Value = 616.625
$ws.Range("J64").Value = 768.9
$ws.Range("K64").Value = 616.625
$ws.Range("L64").Value = 768.9
$ws.Range("M64").Value = -391.625
$ws.Range("N64").Value = -1218.9
$ws.Range("H67").Value = 701.2222
$ws.Range("I67").Value = 616.625
$ws.Range("J67").Value = 768.9
$ws.Range("K67").Value = 616.625
$ws.Range("L67").Value = 768.9
$ws.Range("M67").Value = 163.375
$ws.Range("N67").Value = -2328.9
$ws.Range("H99").Value = 6962.2173
$ws.Range("I99").Value = 7631.55
$ws.Range("K99").Value = 7631.55
$ws.Range("M99").Value = -6133.55
$ws.Range("H134").Value = 1975.1621
$ws.Range("I134").Value = 1745.1724
$ws.Range("K134").Value = 5235.5172
$ws.Range("M134").Value = -2700.5172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9524363
$ws.Range("I16").Value = 14286269
$ws.Range("J16").Value = 552.6
$ws.Range("K16").Value = 14286269
$ws.Range("L16").Value = 552.6
$ws.Range("M16").Value = -14285982
$ws.Range("N16").Value = -1126.6
$ws.Range("H31").Value = 3040.4
$ws.Range("I31").Value = 4591.9
$ws.Range("K31").Value = 4591.9
$ws.Range("M31").Value = -4296.9
$ws.Range("H34").Value = 3040.4
$ws.Range("I34").Value = 4591.9
$ws.Range("K34").Value = 4591.9
$ws.Range("M34").Value = -4389.9
$ws.Range("H113").Value = 9524363
$ws.Range("I113").Value = 14286269
$ws.Range("J113").Value = 552.6
$ws.Range("K113").Value = 14286269
$ws.Range("L113").Value = 552.6
$ws.Range("M113").Value = -14284099
$ws.Range("N113").Value = -4892.6
$ws.Range("H122").Value = 5842.2856
$ws.Range("I122").Value = 5919.4
$ws.Range("K122").Value = 17758.2
$ws.Range("M122").Value = -15308.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1090.091
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3224
$ws.Range("H12").Value = 712
$ws.Range("J12").Value = 881.2857
$ws.Range("L12").Value = 2643.8571
$ws.Range("N12").Value = -2989.8571
$ws.Range("H32").Value = 7749.3
$ws.Range("I32").Value = 2899.8
$ws.Range("J32").Value = 12598.8
$ws.Range("K32").Value = 8699.400000000001
$ws.Range("L32").Value = 37796.39999999999
$ws.Range("M32").Value = -8416.400000000001
$ws.Range("N32").Value = -38362.39999999999
$ws.Range("H113").Value = 250.96
$ws.Range("J113").Value = 237.8
$ws.Range("L113").Value = 713.4000000000001
$ws.Range("N113").Value = -5053.4
$ws.Range("H132").Value = 2659.0967
$ws.Range("J132").Value = 3394.5386
$ws.Range("L132").Value = 30550.8474
$ws.Range("N132").Value = -35610.8474
$ws.Range("H135").Value = 1090.091
$ws.Range("J135").Value = 1000
$ws.Range("L135").Value = 9000
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 9629.091
$ws.Range("I99").Value = 7092.1
$ws.Range("J99").Value = 34999
$ws.Range("K99").Value = 7092.1
$ws.Range("L99").Value = 34999
$ws.Range("M99").Value = -4846.1
$ws.Range("N99").Value = -39491
$ws.Range("H132").Value = 6580910.5
$ws.Range("I132").Value = 2728.9
$ws.Range("K132").Value = 8186.700000000001
$ws.Range("M132").Value = -5656.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 60507.94
$ws.Range("I81").Value = 1789.75
$ws.Range("K81").Value = 3579.5
$ws.Range("M81").Value = -2518.5
$ws.Range("H84").Value = 60507.94
$ws.Range("I84").Value = 1789.75
$ws.Range("K84").Value = 17897.5
$ws.Range("M84").Value = -12593.5
$ws.Range("H132").Value = 1516.1666
$ws.Range("I132").Value = 1465.5807
$ws.Range("K132").Value = 4396.742099999999
$ws.Range("M132").Value = -1866.742099999999
$ws.Range("H136").Value = 32431.719
$ws.Range("I136").Value = 39247.19
$ws.Range("K136").Value = 117741.57
$ws.Range("M136").Value = -115191.57
